$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (fill in first so the "15/06/2021" / "04104016746" shared strings are created
#     before row 2/3 reuse "22/06/2021" / "04104018841 ", matching the target layout) ---
$ws.Range("H9").Value = "'15/06/2021"
$ws.Range("F9").Value = "'04104016746"
$ws.Range("I9").Value = "'08:30"

# --- Row 10 ---
$ws.Range("H10").Value = "'15/06/2021"
$ws.Range("F10").Value = "'04104016746"
$ws.Range("I10").Value = "'08:30"
$ws.Range("T10").Value = "No"
$ws.Range("U10").Value = "Parcial (Resto del Vehiculo)"
$ws.Range("V10").ClearContents()
$ws.Range("Y10").Value = "Sí"
$ws.Range("Z10").Value = "Rueda"
$ws.Range("AA10").Value = 27433

# --- Row 2 ---
$ws.Range("H2").Value = "'22/06/2021"
$ws.Range("F2").Value = "'04104018841 "
$ws.Range("T2").Value = "No"
$ws.Range("U2").Value = "Parcial (Resto del Vehiculo)"
$ws.Range("V2").ClearContents()

# --- Row 3 ---
$ws.Range("H3").Value = "'22/06/2021"
$ws.Range("F3").Value = "'04104018841 "
$ws.Range("U3").Value = "Cerradura"

# --- Selection / view state ---
$ws.Range("F2:I3").Select()
